# Jogos_da_Semana_FlashScore_2024-11-04.xlsx update
#
# The match that used to sit on row 3 (fT8rSK5A, Corinthians x Palmeiras,
# BRAZIL - SERIE A BETANO) dropped off the sheet, so every row below it
# shifts up by one. Deleting that row reproduces the shift (and the new
# dimension A1:BD6) in one shot; everything else is a handful of odds
# that were refreshed for the remaining matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 3 (Corinthians x Palmeiras) - rows 4..7 shift up to 3..6.
$ws.Rows.Item(3).Delete()

# Row 2 (fsUnkChG, Ind. Rivadavia x Rosario Central) - odds refresh, no row shift.
$ws.Range("I2").Value  = 2.8
$ws.Range("N2").Value  = 4.75
$ws.Range("W2").Value  = 5.5
$ws.Range("AA2").Value = 34
$ws.Range("AN2").Value = 4.5
$ws.Range("AX2").Value = 21

# Row 5 (txqKnEdc, Atl. Nacional x Santa Fe, after the shift) - odds refresh.
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63

# Row 6 (rw7N2WTs, Cerro Porteno x Tacuary, after the shift) - odds refresh.
$ws.Range("Q6").Value = 1.88
$ws.Range("R6").Value = 1.93
